# Weekly logic sheet update: a new week of data (row group) is inserted
# right after row 621 and the rest of the data (previously rows 622:747)
# shifts down by one group (3 rows) to rows 625:750.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows before the current row 622, pushing existing
# rows 622:747 down to 625:750 (dimension becomes A1:T750 automatically).
$ws.Rows("622:624").Insert()

# Seed the 3 new rows with the same template as the (now shifted) group
# that used to occupy 622:624 and now sits at 625:627, then overwrite the
# columns that differ for this new week (Fecha, Precio min/max/prom, Origen, Precio $/Kg).
$ws.Range("A625:T627").Copy()
$ws.Range("A622").PasteSpecial()

$ws.Range("D622:D624").Value = 44694
$ws.Range("N622:N624").Value = 7000
$ws.Range("O622:O624").Value = 7500
$ws.Range("P622:P624").Value = 7250
$ws.Range("R622:R624").Value = "Brasil"
$ws.Range("S622:S624").Value = 1812
